$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.719.75'
$ws.Range('E2').Value = '  +1.82%  '
$ws.Range('D3').Value = '1.854.49'
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('D4').Value = "'0.9996"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'244.06"
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').Value = "'0.6396"
$ws.Range('E6').Value = '  +3.39%  '
$ws.Range('D7').Value = "'0.9998"
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = "'46.86"
$ws.Range('E8').Value = '  +3.27%  '
$ws.Range('D9').Value = "'0.2988"
$ws.Range('E9').Value = '  +3.31%  '
$ws.Range('D10').Value = "'0.07476"
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('D11').Value = "'24.28"
$ws.Range('E11').Value = '  +5.11%  '
$ws.Range('D12').Value = "'0.07643"
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '1.871.68'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').Value = "'5.042"
$ws.Range('E14').Value = '  +1.86%  '
$ws.Range('D15').Value = "'0.6873"
$ws.Range('E15').Value = '  +3.36%  '
$ws.Range('D16').Value = "'83.73"
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('D17').Value = "'0.000009493"
$ws.Range('E17').Value = '  +6.60%  '
$ws.Range('D18').Value = "'6.058"
$ws.Range('E18').Value = '  +3.45%  '
$ws.Range('D19').Value = '29.746.37'
$ws.Range('E19').Value = '  +2.04%  '
$ws.Range('D20').Value = '2.118.98'
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('D21').Value = "'235.79"
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D22').Value = "'12.62"
$ws.Range('E22').Value = '  +1.37%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = "'7.400"
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = "'158.16"
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = "'0.1417"
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('D28').Value = "'8.491"
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('D30').Value = "'0.06208"
$ws.Range('E30').Value = '  +5.09%  '
$ws.Range('D31').Value = "'1.494"
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('E32').Value = '  +5.58%  '
$ws.Range('D33').Value = "'4.147"
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('D34').Value = "'4.093"
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('D35').Value = "'1.886"
$ws.Range('E35').Value = '  +1.62%  '
$ws.Range('E36').Value = '  +3.27%  '
$ws.Range('D37').Value = "'0.7282"
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').Value = "'2.607"
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').Value = "'0.01780"
$ws.Range('E40').Value = '  +1.89%  '
$ws.Range('D41').Value = '1.202.46'
$ws.Range('E41').Value = '  -1.55%  '
$ws.Range('D42').Value = "'0.9263"
$ws.Range('E42').Value = '  +0.79%  '
$ws.Range('D43').Value = "'6.145"
$ws.Range('E43').Value = '  -2.19%  '
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = '2.030.59'
$ws.Range('E45').Value = '  +2.89%  '
$ws.Range('D46').Value = "'101.92"
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = "'66.03"
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('D48').Value = "'0.00000000119"
$ws.Range('E48').Value = '  -1.20%  '
$ws.Range('D49').Value = "'0.4057"
$ws.Range('E49').Value = '  +0.94%  '
$ws.Range('D50').Value = "'9.146"
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = "'0.05794"
$ws.Range('E51').Value = '  +0.77%  '
